# Add a new "2022-Q4" quarterly sheet (with its fund-holding detail) and
# record its summary row at the top of the "总计" sheet, pushing the
# existing quarters down by one row/position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q4.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.03

# Style the new index cell (A2) like the rest of the A column.
$total.Range("A2").Font.Bold = $true
$total.Range("A2").Borders.LineStyle = 1
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160

# Renumber the 0-based index column (A) for every data row now that a row
# was inserted at the top.
for ($r = 2; $r -le 10; $r++) {
    $total.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund detail behind that summary row.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$rows = @(
    @("540004", "汇丰晋信2026周期混合",       "1.08", "23.97", "1.66", "0.0179", 5),
    @("009263", "华宝红利精选混合A",          "0.27", "91.60", "2.12", "0.0057", 2),
    @("010841", "华宝红利精选混合C",          "0.18", "91.60", "2.12", "0.0038", 2)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $newSheet.Range("A$r")
    $aCell.Value = $r - 2
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    foreach ($col in @("B", "C", "D", "E", "F", "G")) {
        $idx = [array]::IndexOf(@("B", "C", "D", "E", "F", "G"), $col)
        $textCell = $newSheet.Range("$col$r")
        $textCell.NumberFormat = "@"
        $textCell.Value = $row[$idx]
    }

    $newSheet.Range("H$r").Value = $row[6]

    $r++
}

$newSheet.Range("A1:H" + ($rows.Length + 1)).Columns.AutoFit() | Out-Null
